$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing labels (units added) ---
$ws.Range("B4").Value = "U_cc [V]"
$ws.Range("B5").Value = "R_v [Ohm]"
$ws.Range("B7").Value = "R_ADC [Ohm]"

# --- Move the "U_ADC" formula row from row 9 up to row 8 ---
$ws.Range("B8").Value = "U_ADC [V]"
$ws.Range("C8").NumberFormat = "0.00"
$ws.Range("C8").Formula = "=C4*(C7/(C5+C7))"
$ws.Range("D8").NumberFormat = "0.00"
$ws.Range("D8").Formula = "=C4*(D7/(C5+D7))"
$ws.Range("E8").NumberFormat = "0.00"
$ws.Range("E8").Formula = "=C4*(E7/(C5+E7))"
$ws.Range("F8").NumberFormat = "0.00"
$ws.Range("F8").Formula = "=C4*(F7/(C5+F7))"
$ws.Range("G8").NumberFormat = "0.00"
$ws.Range("G8").Formula = "=C4*(G7/(C5+G7))"
$ws.Range("H8").NumberFormat = "0.00"
$ws.Range("H8").Formula = "=C4*(H7/(C5+H7))"
$ws.Range("I8").NumberFormat = "0.00"
$ws.Range("I8").Formula = "=C4*(I7/(C5+I7))"
$ws.Range("J8").NumberFormat = "0.00"
$ws.Range("J8").Formula = "=C4*(J7/(C5+J7))"

# clear the old row 9 (now vacated)
$ws.Range("B9:J9").Clear()

# --- "ADC Max wert" (was the unlabeled 4096 value in B12) moves to C10, with a new label in B10 ---
$ws.Range("B10").Value = "ADC Max wert"
$ws.Range("C10").Value = 4096

# --- "ADC Step" (was the unlabeled formula in B13) moves to C11, with a new label in B11 ---
$ws.Range("B11").Value = "ADC Step "
$ws.Range("C11").Formula = "=C4/C10"

# clear old B12 / B13
$ws.Range("B12").Clear()
$ws.Range("B13").Clear()

# --- Move the "ADC Wert" formula row from row 15 up to row 13, referencing the new C8/C11 cells ---
$ws.Range("B13").Value = "ADC Wert"
$ws.Range("C13").NumberFormat = "0"
$ws.Range("C13").Formula = "=C8/C11"
$ws.Range("D13").NumberFormat = "0"
$ws.Range("D13").Formula = "=D8/C11"
$ws.Range("E13").NumberFormat = "0"
$ws.Range("E13").Formula = "=E8/C11"
$ws.Range("F13").NumberFormat = "0"
$ws.Range("F13").Formula = "=F8/C11"
$ws.Range("G13").NumberFormat = "0"
$ws.Range("G13").Formula = "=G8/C11"
$ws.Range("H13").NumberFormat = "0"
$ws.Range("H13").Formula = "=H8/C11"
$ws.Range("I13").NumberFormat = "0"
$ws.Range("I13").Formula = "=I8/C11"
$ws.Range("J13").NumberFormat = "0"
$ws.Range("J13").Formula = "=J8/C11"

# clear the old row 15 (now vacated)
$ws.Range("B15:J15").Clear()

# --- Column B now needs an explicit width (grew to fit the new longer labels) ---
$ws.Range("B1").ColumnWidth = 11.3

# --- Page setup: A4, portrait ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Selection moves to C10 ---
$null = $ws.Range("C10").Select()
